$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H64").Value = 4058
$wsALC.Range("I64").Value = 4196
$wsALC.Range("K64").Value = 4196
$wsALC.Range("M64").Value = -3948
$wsALC.Range("H67").Value = 4058
$wsALC.Range("I67").Value = 4196
$wsALC.Range("K67").Value = 4196
$wsALC.Range("M67").Value = -3338
$wsALC.Range("H74").Value = 3237.75
$wsALC.Range("I74").Value = 3000.6
$wsALC.Range("K74").Value = 3000.6
$wsALC.Range("M74").Value = -2064.6
$wsALC.Range("H77").Value = 3237.75
$wsALC.Range("I77").Value = 3000.6
$wsALC.Range("K77").Value = 15003
$wsALC.Range("M77").Value = -10323

$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H2").Value = 1152.5652
$wsARM.Range("I2").Value = 919.9375
$wsARM.Range("K2").Value = 919.9375
$wsARM.Range("M2").Value = -806.9375
$wsARM.Range("H32").Value = 7631.755
$wsARM.Range("I32").Value = 5868.274
$wsARM.Range("J32").Value = 18212.643
$wsARM.Range("K32").Value = 5868.274
$wsARM.Range("L32").Value = 18212.643
$wsARM.Range("M32").Value = -5581.274
$wsARM.Range("N32").Value = -18786.643
$wsARM.Range("H45").Value = 1298.3334
$wsARM.Range("I45").Value = 1286.25
$wsARM.Range("K45").Value = 1286.25
$wsARM.Range("M45").Value = -909.25
$wsARM.Range("H61").Value = 38462548
$wsARM.Range("I61").Value = 41667636
$wsARM.Range("K61").Value = 41667636
$wsARM.Range("M61").Value = -41667424
$wsARM.Range("H116").Value = 1152.5652
$wsARM.Range("I116").Value = 919.9375
$wsARM.Range("K116").Value = 919.9375
$wsARM.Range("M116").Value = 1374.0625
$wsARM.Range("H122").Value = 2280.6316
$wsARM.Range("I122").Value = 1964.9412
$wsARM.Range("K122").Value = 5894.8236
$wsARM.Range("M122").Value = -3444.8236
$wsARM.Range("H136").Value = 38462548
$wsARM.Range("I136").Value = 41667636
$wsARM.Range("K136").Value = 125002908
$wsARM.Range("M136").Value = -125000358

$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("H3").Value = 1152.5652
$wsBSM.Range("I3").Value = 919.9375
$wsBSM.Range("K3").Value = 919.9375
$wsBSM.Range("M3").Value = -805.9375
$wsBSM.Range("H94").Value = 10417027
$wsBSM.Range("I94").Value = 10417027
$wsBSM.Range("K94").Value = 10417027
$wsBSM.Range("M94").Value = -10416576
$wsBSM.Range("H105").Value = 90910550
$wsBSM.Range("I105").Value = 125001110
$wsBSM.Range("J105").Value = 2400
$wsBSM.Range("K105").Value = 125001110
$wsBSM.Range("L105").Value = 2400
$wsBSM.Range("M105").Value = -124999363
$wsBSM.Range("N105").Value = -5894
$wsBSM.Range("H134").Value = 1320.6471
$wsBSM.Range("I134").Value = 1046
$wsBSM.Range("J134").Value = 1979.8
$wsBSM.Range("K134").Value = 3138
$wsBSM.Range("L134").Value = 5939.4
$wsBSM.Range("M134").Value = -603
$wsBSM.Range("N134").Value = -11009.4

$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("H31").Value = 1193.6177
$wsCRP.Range("I31").Value = 1147.0702
$wsCRP.Range("K31").Value = 1147.0702
$wsCRP.Range("M31").Value = -852.0702000000001
$wsCRP.Range("H34").Value = 1193.6177
$wsCRP.Range("I34").Value = 1147.0702
$wsCRP.Range("K34").Value = 1147.0702
$wsCRP.Range("M34").Value = -945.0702000000001
$wsCRP.Range("H58").Value = 5424.9287
$wsCRP.Range("I58").Value = 994.93335
$wsCRP.Range("J58").Value = 10536.462
$wsCRP.Range("K58").Value = 994.93335
$wsCRP.Range("L58").Value = 10536.462
$wsCRP.Range("M58").Value = -791.93335
$wsCRP.Range("N58").Value = -10942.462
$wsCRP.Range("H108").Value = 30995
$wsCRP.Range("J108").Value = 33588.5
$wsCRP.Range("L108").Value = 33588.5
$wsCRP.Range("N108").Value = -41268.5
$wsCRP.Range("H134").Value = 27779526
$wsCRP.Range("I134").Value = 1662.5
$wsCRP.Range("J134").Value = 125002050
$wsCRP.Range("K134").Value = 4987.5
$wsCRP.Range("L134").Value = 375006150
$wsCRP.Range("M134").Value = -2452.5
$wsCRP.Range("N134").Value = -375011220
$wsCRP.Range("H136").Value = 5424.9287
$wsCRP.Range("I136").Value = 994.93335
$wsCRP.Range("J136").Value = 10536.462
$wsCRP.Range("K136").Value = 2984.80005
$wsCRP.Range("L136").Value = 31609.386
$wsCRP.Range("M136").Value = -434.8000499999998
$wsCRP.Range("N136").Value = -36709.386

$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H33").Value = 470.2
$wsCUL.Range("J33").Value = 667
$wsCUL.Range("L33").Value = 4002
$wsCUL.Range("N33").Value = -4568
$wsCUL.Range("H70").Value = 10020.526
$wsCUL.Range("I70").Value = 19300.166
$wsCUL.Range("K70").Value = 57900.49800000001
$wsCUL.Range("M70").Value = -57585.49800000001
$wsCUL.Range("H73").Value = 10020.526
$wsCUL.Range("I73").Value = 19300.166
$wsCUL.Range("K73").Value = 57900.49800000001
$wsCUL.Range("M73").Value = -56808.49800000001
$wsCUL.Range("H131").Value = 23812784
$wsCUL.Range("J131").Value = 4270.7095
$wsCUL.Range("L131").Value = 12812.1285
$wsCUL.Range("N131").Value = -22892.1285

$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("H80").Value = 5497.778
$wsGSM.Range("I80").Value = 5280
$wsGSM.Range("J80").Value = 5770
$wsGSM.Range("K80").Value = 5280
$wsGSM.Range("L80").Value = 5770
$wsGSM.Range("M80").Value = -4282
$wsGSM.Range("N80").Value = -7766
$wsGSM.Range("H83").Value = 5497.778
$wsGSM.Range("I83").Value = 5280
$wsGSM.Range("J83").Value = 5770
$wsGSM.Range("K83").Value = 26400
$wsGSM.Range("L83").Value = 28850
$wsGSM.Range("M83").Value = -21408
$wsGSM.Range("N83").Value = -38834
$wsGSM.Range("H102").Value = 1464.963
$wsGSM.Range("I102").Value = 1517.5652
$wsGSM.Range("K102").Value = 1517.5652
$wsGSM.Range("M102").Value = 104.4348
$wsGSM.Range("H122").Value = 4154.5454
$wsGSM.Range("I122").Value = 4154.5454
$wsGSM.Range("J122").Value = 0
$wsGSM.Range("K122").Value = 12463.6362
$wsGSM.Range("L122").Value = 0
$wsGSM.Range("M122").ClearContents()
$wsGSM.Range("N122").Value = -10013.6362
$wsGSM.Range("H126").Value = 2062
$wsGSM.Range("I126").Value = 1755
$wsGSM.Range("J126").Value = 2325.1428
$wsGSM.Range("K126").Value = 5265
$wsGSM.Range("L126").Value = 6975.428400000001
$wsGSM.Range("M126").Value = -2795
$wsGSM.Range("N126").Value = -11915.4284
$wsGSM.Range("H132").Value = 2127.2075
$wsGSM.Range("I132").Value = 1805.9
$wsGSM.Range("J132").Value = 3115.8462
$wsGSM.Range("K132").Value = 5417.700000000001
$wsGSM.Range("L132").Value = 9347.5386
$wsGSM.Range("M132").Value = -2887.700000000001
$wsGSM.Range("N132").Value = -14407.5386
$wsGSM.Range("H136").Value = 13024.25
$wsGSM.Range("J136").Value = 13024.25
$wsGSM.Range("L136").Value = 39072.75
$wsGSM.Range("N136").Value = -44172.75

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("H7").Value = 2826
$wsLTW.Range("I7").Value = 2592.6667
$wsLTW.Range("J7").Value = 3351
$wsLTW.Range("K7").Value = 2592.6667
$wsLTW.Range("L7").Value = 3351
$wsLTW.Range("M7").Value = -2480.6667
$wsLTW.Range("N7").Value = -3575
$wsLTW.Range("H40").Value = 4920.4614
$wsLTW.Range("I40").Value = 2411.6
$wsLTW.Range("K40").Value = 2411.6
$wsLTW.Range("M40").Value = -2275.6
$wsLTW.Range("H82").Value = 1956.9445
$wsLTW.Range("I82").Value = 1926.3077
$wsLTW.Range("J82").Value = 2036.6
$wsLTW.Range("K82").Value = 1926.3077
$wsLTW.Range("L82").Value = 2036.6
$wsLTW.Range("M82").Value = -1565.3077
$wsLTW.Range("N82").Value = -2758.6
$wsLTW.Range("H85").Value = 1956.9445
$wsLTW.Range("I85").Value = 1926.3077
$wsLTW.Range("J85").Value = 2036.6
$wsLTW.Range("K85").Value = 1926.3077
$wsLTW.Range("L85").Value = 2036.6
$wsLTW.Range("M85").Value = -678.3077000000001
$wsLTW.Range("N85").Value = -4532.6
$wsLTW.Range("H126").Value = 2826
$wsLTW.Range("I126").Value = 2592.6667
$wsLTW.Range("J126").Value = 3351
$wsLTW.Range("K126").Value = 7778.000100000001
$wsLTW.Range("L126").Value = 10053
$wsLTW.Range("M126").Value = -5308.000100000001
$wsLTW.Range("N126").Value = -14993

$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Range("H75").Value = 34500
$wsWVR.Range("J75").Value = 34500
$wsWVR.Range("L75").Value = 34500
$wsWVR.Range("N75").Value = -36372
$wsWVR.Range("H78").Value = 34500
$wsWVR.Range("J78").Value = 34500
$wsWVR.Range("L78").Value = 103500
$wsWVR.Range("N78").Value = -112860
$wsWVR.Range("H126").Value = 166667970
$wsWVR.Range("I126").Value = 333333760
$wsWVR.Range("J126").Value = 2168.3333
$wsWVR.Range("K126").Value = 1000001280
$wsWVR.Range("L126").Value = 6504.999899999999
$wsWVR.Range("M126").Value = -999998810
$wsWVR.Range("N126").Value = -11444.9999
